$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B17").Value = "Hecho"
$ws.Range("B18").Value = "Hecho"

$ws.Application.ActiveWindow.Zoom = 100
[void]$ws.Range("C18").Select()
